$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 5243.6
$ws.Range("I6").Value = 16778.666
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 50335.99800000001
$ws.Range("L6").Value = 900
$ws.Range("M6").Value = -50223.99800000001
$ws.Range("N6").Value = -1124
$ws.Range("H17").Value = 557.413
$ws.Range("J17").Value = 557.413
$ws.Range("L17").Value = 1672.239
$ws.Range("N17").Value = -2008.239
$ws.Range("H33").Value = 171.25
$ws.Range("I33").Value = 106.30769
$ws.Range("J33").Value = 452.66666
$ws.Range("K33").Value = 106.30769
$ws.Range("L33").Value = 452.66666
$ws.Range("M33").Value = 122.69231
$ws.Range("N33").Value = -910.66666
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 3625
$ws.Range("K40").Value = 3625
$ws.Range("M40").Value = -3450
$ws.Range("H113").Value = 25002788
$ws.Range("I113").Value = 66668900
$ws.Range("J113").Value = 3119.8
$ws.Range("K113").Value = 66668900
$ws.Range("L113").Value = 3119.8
$ws.Range("M113").Value = -66665646
$ws.Range("N113").Value = -9627.799999999999
$ws.Range("H116").Value = 3667.5715
$ws.Range("I116").Value = 3322.2222
$ws.Range("J116").Value = 4289.2
$ws.Range("K116").Value = 3322.2222
$ws.Range("L116").Value = 4289.2
$ws.Range("M116").Value = 119.7777999999998
$ws.Range("N116").Value = -11173.2
$ws.Range("H132").Value = 8552110
$ws.Range("I132").Value = 11500303
$ws.Range("K132").Value = 34500909
$ws.Range("M132").Value = -34498379

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3952.5781
$ws.Range("I32").Value = 3642.8728
$ws.Range("J32").Value = 5845.222
$ws.Range("K32").Value = 3642.8728
$ws.Range("L32").Value = 5845.222
$ws.Range("M32").Value = -3355.8728
$ws.Range("N32").Value = -6419.222
$ws.Range("H45").Value = 1111.7037
$ws.Range("I45").Value = 1093.7333
$ws.Range("K45").Value = 1093.7333
$ws.Range("M45").Value = -716.7333000000001
$ws.Range("H61").Value = 1409.1875
$ws.Range("I61").Value = 1262.0769
$ws.Range("J61").Value = 2046.6666
$ws.Range("K61").Value = 1262.0769
$ws.Range("L61").Value = 2046.6666
$ws.Range("M61").Value = -1050.0769
$ws.Range("N61").Value = -2470.6666
$ws.Range("H136").Value = 1409.1875
$ws.Range("I136").Value = 1262.0769
$ws.Range("J136").Value = 2046.6666
$ws.Range("K136").Value = 3786.2307
$ws.Range("L136").Value = 6139.9998
$ws.Range("M136").Value = -1236.2307
$ws.Range("N136").Value = -11239.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 9990
$ws.Range("J18").Value = 9990
$ws.Range("L18").Value = 9990
$ws.Range("N18").Value = -11048

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 1500
$ws.Range("J26").Value = 1500
$ws.Range("L26").Value = 1500
$ws.Range("N26").Value = -2074
$ws.Range("H31").Value = 1271.1964
$ws.Range("I31").Value = 1215.1346
$ws.Range("K31").Value = 1215.1346
$ws.Range("M31").Value = -920.1346000000001
$ws.Range("H34").Value = 1271.1964
$ws.Range("I34").Value = 1215.1346
$ws.Range("K34").Value = 1215.1346
$ws.Range("M34").Value = -1013.1346
$ws.Range("H99").Value = 2082.4
$ws.Range("I99").Value = 2003
$ws.Range("J99").Value = 2400
$ws.Range("K99").Value = 2003
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -505
$ws.Range("N99").Value = -5396
$ws.Range("H126").Value = 2082.4
$ws.Range("I126").Value = 2003
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 6009
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -3539
$ws.Range("N126").Value = -12140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 836246.4
$ws.Range("I4").Value = 179853.8
$ws.Range("J4").Value = 1200908.9
$ws.Range("K4").Value = 539561.3999999999
$ws.Range("L4").Value = 3602726.7
$ws.Range("M4").Value = -539449.3999999999
$ws.Range("N4").Value = -3602950.7
$ws.Range("H86").Value = 300
$ws.Range("I86").Value = 300
$ws.Range("J86").Value = 300
$ws.Range("K86").Value = 900
$ws.Range("L86").Value = 900
$ws.Range("M86").Value = 286
$ws.Range("N86").Value = -3272
$ws.Range("H89").Value = 300
$ws.Range("I89").Value = 300
$ws.Range("J89").Value = 300
$ws.Range("K89").Value = 2700
$ws.Range("L89").Value = 2700
$ws.Range("M89").Value = 3228
$ws.Range("N89").Value = -14556
$ws.Range("H131").Value = 18184672
$ws.Range("I131").Value = 250000370
$ws.Range("J131").Value = 3049.2354
$ws.Range("K131").Value = 750001110
$ws.Range("L131").Value = 9147.706200000001
$ws.Range("M131").Value = -749996070
$ws.Range("N131").Value = -19227.7062

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1179.55
$ws.Range("I113").Value = 1039.6923
$ws.Range("J113").Value = 1439.2858
$ws.Range("K113").Value = 1039.6923
$ws.Range("L113").Value = 1439.2858
$ws.Range("M113").Value = 1130.3077
$ws.Range("N113").Value = -5779.2858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5509
$ws.Range("I46").Value = 2433
$ws.Range("K46").Value = 2433
$ws.Range("M46").Value = -2245
$ws.Range("H55").Value = 238.92592
$ws.Range("I55").Value = 196.73334
$ws.Range("K55").Value = 196.73334
$ws.Range("M55").Value = -23.73334
$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("N94").Value = -21352
$ws.Range("H132").Value = 26874.9
$ws.Range("I132").Value = 1449.7084
$ws.Range("J132").Value = 65012.688
$ws.Range("K132").Value = 4349.1252
$ws.Range("L132").Value = 195038.064
$ws.Range("M132").Value = -1819.1252
$ws.Range("N132").Value = -200098.064

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9286589
$ws.Range("I122").Value = 11819077
$ws.Range("J122").Value = 800.8333
$ws.Range("K122").Value = 35457231
$ws.Range("L122").Value = 2402.4999
$ws.Range("M122").Value = -35454781
$ws.Range("N122").Value = -7302.4999
$ws.Range("H133").Value = 49600
$ws.Range("J133").Value = 49600
$ws.Range("L133").Value = 49600
$ws.Range("N133").Value = -59720
